$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "index" header/column to "i" (updates shared strings + table column name)
$ws.Range("A1").Value = "i"

# Decrement every row index in column A by 1 (was 1-based, now 0-based)
for ($r = 2; $r -le 503; $r++) {
    $c = $ws.Cells.Item($r, 1)
    $c.Value = $c.Value() - 1
}

# Narrow column A now that values/labels are shorter
$ws.Columns.Item(1).ColumnWidth = 3.17
